$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Di aplikasi sistem penunjang keputusan penerimaan calon asisten Laboratorium Teknik Informatika mengimplementasikan algoritma ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Di aplikasi sistem penunjang keputusan penerimaan calon asisten Laboratorium Teknik Informatika mengimplementasikan algoritma ",
    2) | Out-Null

Write-Output "done"
